$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.828.29'
$ws.Range("E2").Value = '  -4.68%  '
$ws.Range("D3").Value = '3.218.18'
$ws.Range("E3").Value = '  -5.78%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '173.96'
$ws.Range("E5").Value = '  -5.55%  '
$ws.Range("D6").Value = '513.49'
$ws.Range("E6").Value = '  -4.25%  '
$ws.Range("E7").Value = '  -4.23%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '3.213.62'
$ws.Range("E9").Value = '  -5.73%  '
$ws.Range("D10").Value = '0.598'
$ws.Range("E10").Value = '  -5.62%  '
$ws.Range("D11").Value = '52.46'
$ws.Range("E11").Value = '  -9.67%  '
$ws.Range("E13").Value = '  -2.81%  '
$ws.Range("E14").Value = '  -6.36%  '
$ws.Range("D15").Value = '3.741.93'
$ws.Range("E15").Value = '  -5.40%  '
$ws.Range("E16").Value = '  -6.56%  '
$ws.Range("D17").Value = '3.221.85'
$ws.Range("E17").Value = '  -5.77%  '
$ws.Range("D18").Value = '62.803.64'
$ws.Range("E18").Value = '  -4.32%  '
$ws.Range("D19").Value = '17.10'
$ws.Range("E19").Value = '  -3.46%  '
$ws.Range("D20").Value = '10.92'
$ws.Range("E20").Value = '  -4.37%  '
$ws.Range("D21").Value = '0.955'
$ws.Range("E21").Value = '  -3.66%  '
$ws.Range("D22").Value = '364.77'
$ws.Range("E22").Value = '  -4.65%  '
$ws.Range("D23").Value = '3.69'
$ws.Range("E23").Value = '  -2.72%  '
$ws.Range("D24").Value = '79.93'
$ws.Range("E24").Value = '  -4.59%  '
$ws.Range("D25").Value = '11.01'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '3.90'
$ws.Range("E26").Value = '  +5.22%  '
$ws.Range("D27").Value = '6.09'
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("D28").Value = '2.60'
$ws.Range("E28").Value = '  -4.46%  '
$ws.Range("D29").Value = '11.17'
$ws.Range("E29").Value = '  -5.69%  '
$ws.Range("D30").Value = '8.13'
$ws.Range("E30").Value = '  -5.73%  '
$ws.Range("D31").Value = '650.62'
$ws.Range("E31").Value = '  -6.84%  '
$ws.Range("D32").Value = '28.09'
$ws.Range("E32").Value = '  -6.45%  '
$ws.Range("D33").Value = '6.28'
$ws.Range("E33").Value = '  -8.84%  '
$ws.Range("D34").Value = '11.07'
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("D36").Value = '57.23'
$ws.Range("E36").Value = '  -7.47%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = '36.49'
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").Value = '0.372'
$ws.Range("E39").Value = '  -5.47%  '
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("E41").Value = '  +11.07%  '
$ws.Range("E42").Value = '  -4.93%  '
$ws.Range("D43").Value = '2.856.57'
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("D44").Value = '2.50'
$ws.Range("E44").Value = '  +3.71%  '
$ws.Range("D45").Value = '2.66'
$ws.Range("E45").Value = '  -1.17%  '
$ws.Range("D46").Value = '2.84'
$ws.Range("E46").Value = '  +8.59%  '
$ws.Range("D47").Value = '0.0389'
$ws.Range("E47").Value = '  -3.43%  '
$ws.Range("E48").Value = '  -8.18%  '
$ws.Range("D49").Value = '135.29'
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("E50").Value = '  -3.39%  '
$ws.Range("D51").Value = '2.92'
$ws.Range("E51").Value = '  +0.63%  '
